$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting existing columns C:X to the right.
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = 14.6

# Set values for new chat_id column (data rows first, then header) so
# the shared-string table grows in the same order as the target file.
$ws.Range("C2").Value = "enemy_general_1"
$ws.Range("C3").Value = "tanooki_1"
$ws.Range("C4").Value = "enemy_general_1"
$ws.Range("C1").Value = "chat_id"

# Update selection to match target workbook state.
$ws.Range("C5").Select()
